$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts) values per row. Regenerate to the
# corrected per-game values (previously mis-derived).
$kValues = @{
    2  = 5
    3  = 11
    4  = 2
    5  = 4
    6  = 4
    7  = 5
    8  = 6
    9  = 5
    10 = 4
    11 = 6
    12 = 2
    13 = 5
    14 = 3
    15 = 1
    16 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
